$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-04 Sunday" "2025-05-05 Monday"

Replace-Text "572×2=1144" "839×5=4195"
Replace-Text "164×4=656" "407×5=2035"
Replace-Text "989×6=5934" "458×7=3206"
Replace-Text "601×7=4207" "379×6=2274"
Replace-Text "851×4=3404" "309×8=2472"

Replace-Text "240×9=2160" "291×6=1746"
Replace-Text "458×6=2748" "907×3=2721"
Replace-Text "562×2=1124" "276×6=1656"
Replace-Text "172×2=344" "490×7=3430"
Replace-Text "180×5=900" "366×9=3294"

Replace-Text "160×4=640" "570×6=3420"
Replace-Text "616×5=3080" "573×9=5157"
Replace-Text "247×8=1976" "317×3=951"
Replace-Text "147×8=1176" "951×7=6657"
Replace-Text "110×8=880" "645×3=1935"

Replace-Text "513×5=2565" "947×3=2841"
Replace-Text "593×8=4744" "613×8=4904"
Replace-Text "905×6=5430" "847×6=5082"
Replace-Text "948×6=5688" "542×5=2710"
Replace-Text "349×8=2792" "152×7=1064"

Replace-Text "680×5=3400" "630×8=5040"
Replace-Text "795×4=3180" "201×8=1608"
Replace-Text "212×7=1484" "510×4=2040"
Replace-Text "493×8=3944" "976×6=5856"
Replace-Text "775×3=2325" "634×8=5072"
